# Generate Report for Handoff
#
# The handback status of the "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md" file
# moved from "Handed back: in sync with en-US" to "Ready for handoff", with
# an updated "Latest Handoff Datetime"/"Latest HO Xliff Generate Date" and a
# new error message explaining that the handback file is stale. This touches
# the Overview sheet plus the per-language detail sheets (zh-cn / de-de).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d47472660c14f89d3bc9c8c0b32d44401898fae/e2e/f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2920032f77682b14ad7dc7e4238cf6dcec41e5b5/e2e/f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md."

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 12:47:11"

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-19 12:47:02"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de detail sheet -----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-19 12:47:11"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
